# Update FFXIV Twintania leve profit market data across all job sheets
# (scheduled runner refresh of currentAveragePrice / profit columns)

$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null
$ws.Range("H33").Value = 846697.3
$ws.Range("I33").Value = 1229390.9
$ws.Range("J33").Value = 4771.4
$ws.Range("K33").Value = 1229390.9
$ws.Range("L33").Value = 4771.4
$ws.Range("M33").Value = -1229161.9
$ws.Range("N33").Value = -5229.4
$ws.Range("H74").Value = 7900
$ws.Range("I74").Value = 7900
$ws.Range("K74").Value = 7900
$ws.Range("M74").Value = -6964
$ws.Range("H77").Value = 7900
$ws.Range("I77").Value = 7900
$ws.Range("K77").Value = 39500
$ws.Range("M77").Value = -34820
$ws.Range("H92").Value = 614
$ws.Range("I92").Value = 637.875
$ws.Range("J92").Value = 518.5
$ws.Range("K92").Value = 637.875
$ws.Range("L92").Value = 518.5
$ws.Range("M92").Value = 610.125
$ws.Range("N92").Value = -3014.5
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null
$ws.Range("H106").Value = 19553.088
$ws.Range("I106").Value = 3826
$ws.Range("K106").Value = 3826
$ws.Range("M106").Value = -3195
$ws.Range("H126").Value = 90000
$ws.Range("J126").Value = 90000
$ws.Range("L126").Value = 90000
$ws.Range("N126").Value = -99880
$ws.Range("H127").Value = 167993.25
$ws.Range("I127").Value = 201093.6
$ws.Range("J127").Value = 2491.5
$ws.Range("K127").Value = 603280.8
$ws.Range("L127").Value = 7474.5
$ws.Range("M127").Value = -598320.8
$ws.Range("N127").Value = -17394.5
$ws.Range("H129").Value = 3361.5715
$ws.Range("I129").Value = 699.4167
$ws.Range("J129").Value = 6911.1113
$ws.Range("K129").Value = 2098.2501
$ws.Range("L129").Value = 20733.3339
$ws.Range("M129").Value = 2901.7499
$ws.Range("N129").Value = -30733.3339
$ws.Range("H131").Value = 12704.571
$ws.Range("I131").Value = 12656
$ws.Range("J131").Value = 12996
$ws.Range("K131").Value = 37968
$ws.Range("L131").Value = 38988
$ws.Range("M131").Value = -32928
$ws.Range("N131").Value = -49068
$ws.Range("H134").Value = 159997
$ws.Range("J134").Value = 159997
$ws.Range("L134").Value = 159997
$ws.Range("N134").Value = -170137
$ws.Range("H135").Value = 1933.0714
$ws.Range("I135").Value = 2062.5557
$ws.Range("K135").Value = 18563.0013
$ws.Range("M135").Value = -16028.0013
$ws.Range("H137").Value = 19161.2
$ws.Range("I137").Value = 8260.416999999999
$ws.Range("K137").Value = 24781.251
$ws.Range("M137").Value = -22231.251
$ws.Range("H140").Value = 500000
$ws.Range("J140").Value = 500000
$ws.Range("L140").Value = 500000
$ws.Range("N140").Value = -510360

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 11191.214
$ws.Range("I45").Value = 12559.833
$ws.Range("J45").Value = 2979.5
$ws.Range("K45").Value = 12559.833
$ws.Range("L45").Value = 2979.5
$ws.Range("M45").Value = -12182.833
$ws.Range("N45").Value = -3733.5
$ws.Range("H74").Value = 7670.8237
$ws.Range("J74").Value = 4107.6665
$ws.Range("L74").Value = 4107.6665
$ws.Range("N74").Value = -5855.6665
$ws.Range("H77").Value = 7670.8237
$ws.Range("J77").Value = 4107.6665
$ws.Range("L77").Value = 20538.3325
$ws.Range("N77").Value = -29274.3325

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7602.62
$ws.Range("I134").Value = 3547.8948
$ws.Range("J134").Value = 20442.584
$ws.Range("K134").Value = 10643.6844
$ws.Range("L134").Value = 61327.75199999999
$ws.Range("M134").Value = -8108.6844
$ws.Range("N134").Value = -66397.75199999999

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 44997.5
$ws.Range("J50").Value = 44997.5
$ws.Range("L50").Value = 44997.5
$ws.Range("N50").Value = -46247.5
$ws.Range("H59").Value = 43894.375
$ws.Range("I59").Value = 37299.332
$ws.Range("J59").Value = 45416.31
$ws.Range("K59").Value = 37299.332
$ws.Range("L59").Value = 45416.31
$ws.Range("M59").Value = -36154.332
$ws.Range("N59").Value = -47706.31
$ws.Range("H74").Value = 66661.664
$ws.Range("J74").Value = 66661.664
$ws.Range("L74").Value = 66661.664
$ws.Range("N74").Value = -68409.664
$ws.Range("H77").Value = 66661.664
$ws.Range("J77").Value = 66661.664
$ws.Range("L77").Value = 199984.992
$ws.Range("N77").Value = -208720.992
$ws.Range("H95").Value = 5240
$ws.Range("J95").Value = 5377.6665
$ws.Range("L95").Value = 5377.6665
$ws.Range("N95").Value = -10869.6665

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 258.75
$ws.Range("I18").Value = 276.7143
$ws.Range("K18").Value = 830.1428999999999
$ws.Range("M18").Value = -661.1428999999999
$ws.Range("H94").Value = 6286.875
$ws.Range("I94").Value = 7120
$ws.Range("J94").Value = 5908.1816
$ws.Range("K94").Value = 21360
$ws.Range("L94").Value = 17724.5448
$ws.Range("M94").Value = -20684
$ws.Range("N94").Value = -19076.5448
$ws.Range("H96").Value = 12500
$ws.Range("J96").Value = 12500
$ws.Range("L96").Value = 37500
$ws.Range("N96").Value = -41618
$ws.Range("H100").Value = 9400.556
$ws.Range("I100").Value = 5303
$ws.Range("J100").Value = 10571.286
$ws.Range("K100").Value = 15909
$ws.Range("L100").Value = 31713.858
$ws.Range("M100").Value = -15098
$ws.Range("N100").Value = -33335.858
$ws.Range("H104").Value = 6297.4287
$ws.Range("I104").Value = 4716.4
$ws.Range("J104").Value = 10250
$ws.Range("K104").Value = 14149.2
$ws.Range("L104").Value = 30750
$ws.Range("M104").Value = -11528.2
$ws.Range("N104").Value = -35992
$ws.Range("H106").Value = 12250
$ws.Range("J106").Value = 12250
$ws.Range("L106").Value = 36750
$ws.Range("N106").Value = -38642
$ws.Range("H112").Value = 1490
$ws.Range("J112").Value = 1490
$ws.Range("L112").Value = 4470
$ws.Range("N112").Value = -6686
$ws.Range("H115").Value = 3652
$ws.Range("I115").Value = 4478
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 13434
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = -12259
$ws.Range("N115").Value = -8350
$ws.Range("H117").Value = 1129.4
$ws.Range("I117").Value = 1178
$ws.Range("J117").Value = 935
$ws.Range("K117").Value = 3534
$ws.Range("L117").Value = 2805
$ws.Range("M117").Value = -92
$ws.Range("N117").Value = -9689
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").Value = $null
$ws.Range("H125").Value = 4999
$ws.Range("I125").Value = 4999
$ws.Range("K125").Value = 14997
$ws.Range("M125").Value = -10077
$ws.Range("H140").Value = 5842.5
$ws.Range("I140").Value = 3361.875
$ws.Range("J140").Value = 9150
$ws.Range("K140").Value = 10085.625
$ws.Range("L140").Value = 27450
$ws.Range("M140").Value = -4905.625
$ws.Range("N140").Value = -37810

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10286.875
$ws.Range("I70").Value = 6898.5
$ws.Range("K70").Value = 6898.5
$ws.Range("M70").Value = -6628.5
$ws.Range("H73").Value = 10286.875
$ws.Range("I73").Value = 6898.5
$ws.Range("K73").Value = 6898.5
$ws.Range("M73").Value = -5962.5
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = $null
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = $null

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1819.0625
$ws.Range("I82").Value = 892.5
$ws.Range("J82").Value = 2745.625
$ws.Range("K82").Value = 892.5
$ws.Range("L82").Value = 2745.625
$ws.Range("M82").Value = -531.5
$ws.Range("N82").Value = -3467.625
$ws.Range("H85").Value = 1819.0625
$ws.Range("I85").Value = 892.5
$ws.Range("J85").Value = 2745.625
$ws.Range("K85").Value = 892.5
$ws.Range("L85").Value = 2745.625
$ws.Range("M85").Value = 355.5
$ws.Range("N85").Value = -5241.625
$ws.Range("H93").Value = 3742.25
$ws.Range("I93").Value = 4103.1113
$ws.Range("J93").Value = 494.5
$ws.Range("K93").Value = 4103.1113
$ws.Range("L93").Value = 494.5
$ws.Range("M93").Value = -2855.1113
$ws.Range("N93").Value = -2990.5

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4835367
$ws.Range("I107").Value = 3103.6667
$ws.Range("J107").Value = 13895861
$ws.Range("K107").Value = 9311.000100000001
$ws.Range("L107").Value = 41687583
$ws.Range("M107").Value = -7391.000100000001
$ws.Range("N107").Value = -41691423
$ws.Range("H136").Value = 1518.5428
$ws.Range("I136").Value = 677.0755
$ws.Range("J136").Value = 4141.9414
$ws.Range("K136").Value = 2031.2265
$ws.Range("L136").Value = 12425.8242
$ws.Range("M136").Value = 518.7734999999998
$ws.Range("N136").Value = -17525.8242
